$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

$ws.Range("E2").Value = 91
$ws.Range("E10").Value = 468
$ws.Range("F10").Value = 229
$ws.Range("H10").Value = 229
$ws.Range("E11").Value = 317
$ws.Range("E12").Value = 456
$ws.Range("F12").Value = 249
$ws.Range("H12").Value = 249
$ws.Range("E13").Value = 118
$ws.Range("E14").Value = 117
$ws.Range("F14").Value = 60
$ws.Range("H14").Value = 60
$ws.Range("E15").Value = 150
$ws.Range("F15").Value = 65
$ws.Range("H15").Value = 65
$ws.Range("E17").Value = 88
$ws.Range("E20").Value = 82
$ws.Range("E23").Value = 187
$ws.Range("E24").Value = 197
$ws.Range("F24").Value = 106
$ws.Range("H24").Value = 106
$ws.Range("E25").Value = 247
$ws.Range("E26").Value = 148
$ws.Range("F27").Value = 148
$ws.Range("H27").Value = 148
$ws.Range("E28").Value = 188
$ws.Range("F28").Value = 68
$ws.Range("H28").Value = 68
$ws.Range("E29").Value = 160
$ws.Range("E30").Value = 195
$ws.Range("E33").Value = 268
$ws.Range("E34").Value = 201
$ws.Range("F34").Value = 129
$ws.Range("H34").Value = 129
$ws.Range("E35").Value = 134
$ws.Range("E37").Value = 148
$ws.Range("E38").Value = 85
$ws.Range("E40").Value = 243
$ws.Range("E41").Value = 374
$ws.Range("E42").Value = 342
$ws.Range("F42").Value = 184
$ws.Range("H42").Value = 184
$ws.Range("E44").Value = 294
$ws.Range("F44").Value = 146
$ws.Range("H44").Value = 146
$ws.Range("E45").Value = 131
$ws.Range("F45").Value = 65
$ws.Range("H45").Value = 65
$ws.Range("E46").Value = 291
$ws.Range("E47").Value = 419
$ws.Range("E48").Value = 190
$ws.Range("E49").Value = 272
$ws.Range("F49").Value = 116
$ws.Range("H49").Value = 116
$ws.Range("E50").Value = 230
$ws.Range("E51").Value = 218
